$wb = $excel.ActiveWorkbook

# --- Rename sheet "son" -> "userlive" ---
$wsLive = $wb.Worksheets.Item("son")
$wsLive.Name = "userlive"

$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# --- Update Sheet1 (A1:I2) ---
$wsSheet1.Range("C2").Value = "hh_mk_doncap"
$wsSheet1.Range("D2").Value = "hh_ck_doncap"
$wsSheet1.Range("E2").Value = "'016704070003846"
$wsSheet1.Range("G2").Value = "QWJjMTIz"
$wsSheet1.Range("H2").Value = "UXdlMTIz"

$wsSheet1.Range("C2:H2").Select()

# --- Update userlive (A1:H2) ---
$wsLive.Range("C2").Value = "hh_mk_doncap"
$wsLive.Range("D2").Value = "hh_ck_doncap"
$wsLive.Range("E2").Value = "'016704070003846"
$wsLive.Range("G2").Value = "QWJjMTIz"
$wsLive.Range("H2").Value = "UXdlMTIz"

$wsLive.Columns.Item(6).ColumnWidth = 24.736979166666668
$wsLive.Columns.Item(7).ColumnWidth = 14.736979166666666

$wsLive.Range("C2:H2").Select()

$wsSheet1.Activate()
